{"js": "// Ordered list of [oldText, newText] pairs corresponding to each paragraph\n// in document order (body.paragraphs flattens into the title paragraph\n// followed by every table cell's paragraph, in reading order).\nconst replacements = [[\"2023-06-06 Tuesday\", \"2023-06-07 Wednesday\"], [\"51+18=69\", \"75-39=36\"], [\"81-6=75\", \"68-5=63\"], [\"93-11=82\", \"63-18=45\"], [\"34+21=55\", \"76-5=71\"], [\"45+14=59\", \"40+11=51\"], [\"9-8=1\", \"78-9=69\"], [\"27+9=36\", \"31+59=90\"], [\"13+6=19\", \"35-24=11\"], [\"52-39=13\", \"94-67=27\"], [\"59-57=2\", \"64-23=41\"], [\"35+42=77\", \"36+1=37\"], [\"96-65=31\", \"51-5=46\"], [\"84-51=33\", \"1+6=7\"], [\"91-42=49\", \"8+70=78\"], [\"94+5=99\", \"23+23=46\"], [\"11+67=78\", \"41+22=63\"], [\"48-40=8\", \"17-6=11\"], [\"60+24=84\", \"46-38=8\"], [\"23+1=24\", \"3+49=52\"], [\"75-36=39\", \"65-61=4\"], [\"65+21=86\", \"97-13=84\"], [\"37-34=3\", \"33+5=38\"], [\"44+15=59\", \"88-46=42\"], [\"5-4=1\", \"28+12=40\"], [\"40+35=75\", \"47-47=0\"], [\"72-42=30\", \"14+57=71\"], [\"23+2=25\", \"19+41=60\"], [\"66+28=94\", \"61-10=51\"], [\"68+2=70\", \"61-1=60\"], [\"86-69=17\", \"87-8=79\"], [\"7+74=81\", \"11+81=92\"], [\"77-14=63\", \"15+64=79\"], [\"32+4=36\", \"20+31=51\"], [\"51+22=73\", \"93-84=9\"], [\"77-50=27\", \"23-0=23\"], [\"16+33=49\", \"57-56=1\"], [\"34+36=70\", \"2+4=6\"], [\"5+47=52\", \"15+51=66\"], [\"38+18=56\", \"68-1=67\"], [\"71+22=93\", \"99-98=1\"], [\"78-14=64\", \"43-28=15\"], [\"56+12=68\", \"0+20=20\"], [\"95-88=7\", \"65-20=45\"], [\"46+52=98\", \"26+57=83\"], [\"13+73=86\", \"38+32=70\"], [\"46+27=73\", \"38-5=33\"], [\"88-78=10\", \"60-11=49\"], [\"25+18=43\", \"6+17=23\"], [\"24+14=38\", \"4+62=66\"], [\"89-48=41\", \"5+85=90\"], [\"28+44=72\", \"22+6=28\"], [\"2+40=42\", \"23+0=23\"], [\"77-76=1\", \"68+21=89\"], [\"40+14=54\", \"49-18=31\"], [\"76+14=90\", \"65+8=73\"], [\"81-64=17\", \"39-4=35\"], [\"82-61=21\", \"67-53=14\"], [\"15+84=99\", \"84-82=2\"], [\"71+0=71\", \"65-25=40\"], [\"65-54=11\", \"42+50=92\"], [\"36-17=19\", \"29-17=12\"], [\"33-7=26\", \"29+5=34\"], [\"10+43=53\", \"98-77=21\"], [\"33+50=83\", \"29+46=75\"], [\"11+57=68\", \"69-45=24\"], [\"56+1=57\", \"88-76=12\"], [\"12+36=48\", \"28+34=62\"], [\"46+34=80\", \"37+48=85\"], [\"48-28=20\", \"32-25=7\"], [\"92-44=48\", \"94-57=37\"], [\"98-20=78\", \"22-15=7\"], [\"79-29=50\", \"24+29=53\"], [\"36-15=21\", \"64-37=27\"], [\"33-7=26\", \"86-53=33\"], [\"43+51=94\", \"69-49=20\"], [\"12+33=45\", \"99-20=79\"], [\"92-6=86\", \"8+2=10\"], [\"68-15=53\", \"73-29=44\"], [\"4+64=68\", \"98-88=10\"], [\"51-44=7\", \"88-50=38\"], [\"45-42=3\", \"3+78=81\"], [\"73-35=38\", \"13+25=38\"], [\"28+32=60\", \"30-5=25\"], [\"5+46=51\", \"21+44=65\"], [\"45+37=82\", \"80-16=64\"], [\"72-61=11\", \"87-77=10\"], [\"33+9=42\", \"6+65=71\"], [\"69+21=90\", \"24+58=82\"], [\"86-57=29\", \"18+37=55\"], [\"40-8=32\", \"38+27=65\"], [\"69-62=7\", \"52+24=76\"], [\"48-13=35\", \"97-61=36\"], [\"29-21=8\", \"56+8=64\"], [\"30-9=21\", \"67-14=53\"], [\"78-57=21\", \"35+7=42\"], [\"75+3=78\", \"34+43=77\"], [\"17-1=16\", \"86-0=86\"], [\"7+78=85\", \"83-56=27\"], [\"39+35=74\", \"55+3=58\"], [\"31+21=52\", \"58+16=74\"]];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== replacements.length) {\n  throw new Error(\n    `Expected ${replacements.length} paragraphs but found ${items.length}`\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const actual = items[i].text.replace(/[\\r\\n]+$/, \"\");\n  if (actual !== oldText) {\n    throw new Error(\n      `Paragraph ${i} text mismatch: expected \"${oldText}\" but found \"${actual}\"`\n    );\n  }\n  items[i].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$titleOld = '2023-06-06 Tuesday'\n$titleNew = '2023-06-07 Wednesday'\n\n# Each entry: old expression text, new expression text.\n# Cells are in row-major order matching the table's 20 rows x 5 columns.\n$cellPairs = @(\n    @('51+18=69', '75-39=36'),\n    @('81-6=75', '68-5=63'),\n    @('93-11=82', '63-18=45'),\n    @('34+21=55', '76-5=71'),\n    @('45+14=59', '40+11=51'),\n    @('9-8=1', '78-9=69'),\n    @('27+9=36', '31+59=90'),\n    @('13+6=19', '35-24=11'),\n    @('52-39=13', '94-67=27'),\n    @('59-57=2', '64-23=41'),\n    @('35+42=77', '36+1=37'),\n    @('96-65=31', '51-5=46'),\n    @('84-51=33', '1+6=7'),\n    @('91-42=49', '8+70=78'),\n    @('94+5=99', '23+23=46'),\n    @('11+67=78', '41+22=63'),\n    @('48-40=8', '17-6=11'),\n    @('60+24=84', '46-38=8'),\n    @('23+1=24', '3+49=52'),\n    @('75-36=39', '65-61=4'),\n    @('65+21=86', '97-13=84'),\n    @('37-34=3', '33+5=38'),\n    @('44+15=59', '88-46=42'),\n    @('5-4=1', '28+12=40'),\n    @('40+35=75', '47-47=0'),\n    @('72-42=30', '14+57=71'),\n    @('23+2=25', '19+41=60'),\n    @('66+28=94', '61-10=51'),\n    @('68+2=70', '61-1=60'),\n    @('86-69=17', '87-8=79'),\n    @('7+74=81', '11+81=92'),\n    @('77-14=63', '15+64=79'),\n    @('32+4=36', '20+31=51'),\n    @('51+22=73', '93-84=9'),\n    @('77-50=27', '23-0=23'),\n    @('16+33=49', '57-56=1'),\n    @('34+36=70', '2+4=6'),\n    @('5+47=52', '15+51=66'),\n    @('38+18=56', '68-1=67'),\n    @('71+22=93', '99-98=1'),\n    @('78-14=64', '43-28=15'),\n    @('56+12=68', '0+20=20'),\n    @('95-88=7', '65-20=45'),\n    @('46+52=98', '26+57=83'),\n    @('13+73=86', '38+32=70'),\n    @('46+27=73', '38-5=33'),\n    @('88-78=10', '60-11=49'),\n    @('25+18=43', '6+17=23'),\n    @('24+14=38', '4+62=66'),\n    @('89-48=41', '5+85=90'),\n    @('28+44=72', '22+6=28'),\n    @('2+40=42', '23+0=23'),\n    @('77-76=1', '68+21=89'),\n    @('40+14=54', '49-18=31'),\n    @('76+14=90', '65+8=73'),\n    @('81-64=17', '39-4=35'),\n    @('82-61=21', '67-53=14'),\n    @('15+84=99', '84-82=2'),\n    @('71+0=71', '65-25=40'),\n    @('65-54=11', '42+50=92'),\n    @('36-17=19', '29-17=12'),\n    @('33-7=26', '29+5=34'),\n    @('10+43=53', '98-77=21'),\n    @('33+50=83', '29+46=75'),\n    @('11+57=68', '69-45=24'),\n    @('56+1=57', '88-76=12'),\n    @('12+36=48', '28+34=62'),\n    @('46+34=80', '37+48=85'),\n    @('48-28=20', '32-25=7'),\n    @('92-44=48', '94-57=37'),\n    @('98-20=78', '22-15=7'),\n    @('79-29=50', '24+29=53'),\n    @('36-15=21', '64-37=27'),\n    @('33-7=26', '86-53=33'),\n    @('43+51=94', '69-49=20'),\n    @('12+33=45', '99-20=79'),\n    @('92-6=86', '8+2=10'),\n    @('68-15=53', '73-29=44'),\n    @('4+64=68', '98-88=10'),\n    @('51-44=7', '88-50=38'),\n    @('45-42=3', '3+78=81'),\n    @('73-35=38', '13+25=38'),\n    @('28+32=60', '30-5=25'),\n    @('5+46=51', '21+44=65'),\n    @('45+37=82', '80-16=64'),\n    @('72-61=11', '87-77=10'),\n    @('33+9=42', '6+65=71'),\n    @('69+21=90', '24+58=82'),\n    @('86-57=29', '18+37=55'),\n    @('40-8=32', '38+27=65'),\n    @('69-62=7', '52+24=76'),\n    @('48-13=35', '97-61=36'),\n    @('29-21=8', '56+8=64'),\n    @('30-9=21', '67-14=53'),\n    @('78-57=21', '35+7=42'),\n    @('75+3=78', '34+43=77'),\n    @('17-1=16', '86-0=86'),\n    @('7+78=85', '83-56=27'),\n    @('39+35=74', '55+3=58'),\n    @('31+21=52', '58+16=74'),\n)\n\n$d = $word.ActiveDocument\n\n# --- Update the title/date paragraph -----------------------------------\n$titlePara = $d.Paragraphs.Item(1)\n$titleActual = $titlePara.Range.Text.TrimEnd([char]13, [char]10)\nif ($titleActual -ne $titleOld) {\n    throw \"Title paragraph text mismatch: expected '$titleOld' but found '$titleActual'\"\n}\n$titlePara.Range.Text = $titleNew\n\n# --- Update every cell of the first (and only) table --------------------\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\nif (($rowCount * $colCount) -ne $cellPairs.Count) {\n    throw \"Expected $($cellPairs.Count) cells but table has $($rowCount * $colCount)\"\n}\n\n$index = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $pair = $cellPairs[$index]\n        $oldText = $pair[0]\n        $newText = $pair[1]\n\n        $cell = $table.Cell($r, $c)\n        $actual = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($actual -ne $oldText) {\n            throw \"Cell ($r,$c) text mismatch: expected '$oldText' but found '$actual'\"\n        }\n        $cell.Range.Text = $newText\n\n        $index = $index + 1\n    }\n}\n"}
